$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the surviving control-point rows
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 9
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 4

# Remove the trailing rows (4 and 5) entirely
$ws.Range("A4:B5").Delete()
